$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply edits from the bottom of the sheet upward so that earlier row numbers used below
# remain valid (unaffected by the single row insert and two-row delete performed higher up).

# Row 127 -> copyright year line text refresh
$ws.Range("A127").Value = 'li: State Government of Victoria, Australia © 2019'

# Row 125 -> "Last Update" date changes to 21 October 2020
$ws.Range("A125").Value = 'li: Last Update: 21 October 2020'

# Row 118 absorbs the combined bus cancellations text; rows 119-120 (Paschendale Merino
# Bus Route + Term 4 blurb) are removed entirely now that there is no bus cancellation
$ws.Range("A118").Value = 'Bus service cancellations or alterationsThe Department hasnotbeen advised of any bus route cancellations.For Term 4 2020, schools bus services will continue to be provided to support student travel to schools where needed.'
$ws.Rows.Item(119).Resize(2).Delete()

# Row 117: remove stray non-breaking space from TAFE closures text
$ws.Range("A117").Value = 'li: The Department has not been advised of any TAFE closures.'

# Row 115: remove stray non-breaking spaces from "school closures" text
$ws.Range("A115").Value = 'li: been advised of any school closures.'

# Row 80: remove stray non-breaking spaces from early childhood service closures text
$ws.Range("A80").Value = 'li: The Department has been advised of the following early childhood service closures:'

# Row 78: remove stray non-breaking space from "TAFE closures" text
$ws.Range("A78").Value = 'li: been advised of any TAFE closures.'

# Rows 72-74 are replaced with 3 new closure lines (school closures intro + two new school names)
$ws.Range("A74").Value = 'li: Dallas Brooks Community Primary School, DALLAS'
$ws.Range("A73").Value = 'li: East Preston Islamic College, EAST PRESTON'
$ws.Range("A72").Value = 'li: The Department has been advised of the following school closures:'

# Insert a new row before row 61 for the new Dallas Kindergarten closure entry
$ws.Rows.Item(61).Insert()
$ws.Range("A61").Value = 'li: Dallas Kindergarten, DALLAS'

# Row 7: update the "On this page..." banner text with the new date/time
$ws.Range("A7").Value = 'On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school and early childhood service, TAFE closures and relocations for Thursday 22 October, (as at 11.59pm, 22 October)South-Eastern Victoria RegionEarly childhood services'

$ws.Range("A1").Select()
